# Apply the crypto price/volume refresh from the GitHub Actions update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '51.666.13'
$ws.Range('E2').Value = '  +3.79%  '

# Row 3
$ws.Range('D3').Value = '2.751.02'
$ws.Range('E3').Value = '  +3.17%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.09%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '116.96'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +3.43%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '333.90'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.55%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.536'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +1.81%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.08%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.579'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +5.19%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.61'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +2.66%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.15'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.53%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0830'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.27%  '

# Row 13
$ws.Range('E13').Value = '  +2.63%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.63'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +4.13%  '

# Row 15
$ws.Range('D15').Value = '3.176.76'
$ws.Range('E15').Value = '  +2.83%  '

# Row 16
$ws.Range('D16').Value = '2.748.34'
$ws.Range('E16').Value = '  +3.19%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.890'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +2.29%  '

# Row 18
$ws.Range('D18').Value = '51.575.25'
$ws.Range('E18').Value = '  +3.57%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.81'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +5.64%  '

# Row 20
$ws.Range('E20').Value = '  +2.90%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.86'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.69%  '

# Row 22
$ws.Range('D22').Value = '0.0₃0964'
$ws.Range('E22').Value = '  +0.85%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '277.49'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.91%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.36'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.75%  '

# Row 25
$ws.Range('E25').Value = '  +4.38%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.94'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.70%  '

# Row 27
$ws.Range('E27').Value = '  +0.34%  '

# Row 28
$ws.Range('E28').Value = '  +0.16%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.32'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.68%  '

# Row 30
$ws.Range('E30').Value = '  -0.97%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '35.60'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.88%  '

# Row 32
$ws.Range('E32').Value = '  +1.88%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '50.45'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.58%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.63'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +2.89%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0825'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.65%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '19.42'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.06%  '

# Row 38
$ws.Range('E38').Value = '  +2.59%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.31'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +5.82%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.00'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.09%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '129.67'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +3.48%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '23.71'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +6.06%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0347'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +9.56%  '

# Row 44
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.114'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.36%  '

# Row 45
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.29'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +3.75%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.35'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +12.38%  '

# Row 47
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.39'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +3.12%  '

# Row 48
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '2.102.92'
$ws.Range('E48').Value = '  -0.52%  '

# Row 49
$ws.Range('E49').Value = '  +2.56%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.64'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +6.21%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.00'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.02%  '
